# Session9_Assignment_5.docx edit script
# - Inserts a new header row at the top of the table with column titles,
#   including a new last column describing lookup/insert/delete complexity.
# - Fills in the (previously empty) last column for every data row with the
#   Big-O complexity notes.
# - Resizes four of the table columns to their new widths.
# - Tweaks a couple of sentences in the paragraphs that follow the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Insert the new header row above the existing first row -------------
$firstRow = $t.Rows.Item(1)
$newRow = $t.Rows.Add($firstRow)

$t.Cell(1,1).Range.Text = "Data Structure"
$t.Cell(1,2).Range.Text = "Ordering"
$t.Cell(1,3).Range.Text = "Synchronized?"
$t.Cell(1,4).Range.Text = "Duplicates allowed?"
$t.Cell(1,5).Range.Text = "Null allowed?"
$t.Cell(1,6).Range.Text = "Complexity of lookup, insertion,deletion"

# --- 2. Fill the last ("complexity") column for each data row --------------
# Row 2 = ArrayList
$t.Cell(2,6).Range.Text = "O(1) get,add O(n) remove"

# Row 3 = LinkedList (two paragraphs in the cell)
$t.Cell(3,6).Range.Text = "O(n) get,remove" + [char]13 + "O(1) amortized"

# Row 4 = HashSet
$t.Cell(4,6).Range.Text = "O(1) add,remove"

# Row 5 = LinkedHashSet
$t.Cell(5,6).Range.Text = "O(1) add,remove"

# Row 6 = TreeSet
$t.Cell(6,6).Range.Text = "O(log(n)) add,remove"

# Row 7 = HashMap
$t.Cell(7,6).Range.Text = "O(1) get,put,search"

# Row 8 = LinkedHashMap
$t.Cell(8,6).Range.Text = "O(1) get,put"

# Row 9 = TreeMap
$t.Cell(9,6).Range.Text = "O(log(n) get,put"

# --- 3. Resize the columns whose widths changed -----------------------------
$t.Columns.Item(2).Width = 70.05   # 1401 dxa
$t.Columns.Item(4).Width = 71.55   # 1431 dxa
$t.Columns.Item(5).Width = 70.2    # 1404 dxa
$t.Columns.Item(6).Width = 90.3    # 1806 dxa

# --- 4. Prefix the "TreeMap and TreeSet..." note paragraph with "Note:" ----
# (Done via Find/Replace rather than the Paragraphs collection: after the
#  structural table edit above, iterating $d.Paragraphs directly is unreliable.)
$d.Content.Find.Execute(
    "TreeMap and TreeSet enforce",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Note:TreeMap and TreeSet enforce",
    2) | Out-Null

# --- 5. Tweak the "It is not very easy to sort hashmap..." paragraph -------
$d.Content.Find.Execute(
    "no method within hashmap.  We can",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "no method within hashmap or map to do so, like there is for lists and sets.  We can",
    2) | Out-Null

$d.Content.Find.Execute(
    "there is no real way to do so. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "there is no real way to do so using solely a hashmap.",
    2) | Out-Null
